$d = $word.ActiveDocument

# The paragraph currently holds several runs (joined by <w:br/> line breaks):
#   "Hey this is akahara" / "dfdfd" / "fdfdf" / "fdf" /
#   "for jnknknknmnm,nkjkjnjkm,  nkll,kjnm. Asakl.mklmlsa"
# Replace that whole span with a single run of text "jjjjj", leaving the
# trailing _GoBack bookmark untouched.
$d.Content.Find.Execute(
    "Hey this is akahara*Asakl.mklmlsa",  # wildcard match across all runs/line breaks
    $true,   # MatchCase
    $false,  # MatchWholeWord
    $true,   # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap (wdFindContinue)
    $false,  # Format
    "jjjjj", # ReplaceWith
    2        # Replace (wdReplaceAll)
) | Out-Null

Write-Output $d.Content.Text
